# Adds Lesson 16 (3rd ed.) kanji vocabulary rows 54-122 to the wordlist sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    ,('child', '子供|こども')
    ,('to offer something to a spirit', '供える|そなえる')
    ,('offer', '提供|ていきょう')
    ,('the world', '世界|せかい')
    ,('care', '世話|せわ')
    ,('generation', '世代|せだい')
    ,('the third generation', '三世|さんせい')
    ,('society', '世の中|よのなか')
    ,('visibility', '視界|しかい')
    ,('political world', '政界|せいかい')
    ,('limit', '限界|げんかい')
    ,('all; entire', '全部|ぜんぶ')
    ,('safety', '安全|あんぜん')
    ,('whole country', '全国|ぜんこく')
    ,('entirely', '全く|まったく')
    ,('everything; all', '全て|すべて')
    ,('all; entire', '全部|ぜんぶ')
    ,('room', '部屋|へや')
    ,('tennis club', 'テニス部|テニスぶ')
    ,('department manager', '部長|ぶちょう')
    ,('(something) begins', '始まる|はじまる')
    ,('to begin (something)', '始める|はじめる')
    ,('first train (of the day)', '始発|しはつ')
    ,('start', '開始|かいし')
    ,('every week', '毎週|まいしゅう')
    ,('last week', '先週|せんしゅう')
    ,('one week', '一週間|いっしゅうかん')
    ,('second week', '二週目|にしゅうめ')
    ,('weekend', '週末|しゅうまつ')
    ,('to think', '考える|かんがえる')
    ,('idea', '考え|かんがえ')
    ,('archeology', '考古学|こうこがく')
    ,('reference', '参考|さんこう')
    ,('to open (something)', '開ける|あける')
    ,('(something) opens', '開く|あく')
    ,('to open', '開く|ひらく')
    ,('opening of a store', '開店|かいてん')
    ,('room', '部屋|へや')
    ,('bookstore', '本屋|ほんや')
    ,('fish shop', '魚屋|さかなや')
    ,('rooftop', '屋上|おくじょう')
    ,('indoor', '屋内|おくない')
    ,('ally; person on one''s side', '味方|みかた')
    ,('way of reading', '読み方|よみかた')
    ,('evening', '夕方|ゆうがた')
    ,('both', '両方|りょうほう')
    ,('method', '方法|ほうほう')
    ,('exercise', '運動|うんどう')
    ,('driving', '運転|うんてん')
    ,('lucky', '運がいい|うんがいい')
    ,('fate', '運命|うんめい')
    ,('to carry', '運ぶ|はこぶ')
    ,('to move', '動く|うごく')
    ,('automobile', '自動車|じどうしゃ')
    ,('animal', '動物|どうぶつ')
    ,('verb', '動詞|どうし')
    ,('to teach', '教える|おしえる')
    ,('classroom', '教室|きょうしつ')
    ,('church', '教会|きょうかい')
    ,('Christianity', 'キリスト教|キリストきょう')
    ,('textbook', '教科書|きょうかしょ')
    ,('professor''s office', '研究室|けんきゅうしつ')
    ,('basement', '地下室|ちかしつ')
    ,('waiting room', '待合室|まちあいしつ')
    ,('other than...', '～以外|～いがい')
    ,('...or more', '～以上|～いじょう')
    ,('...or less', '～以下|～いか')
    ,('within...', '～以内|～いない')
    ,('before; formerly', '以前|いぜん')
)

$startRow = 54
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $pair = $newRows[$i]
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
}
